# Fix issues resulting from extra files: every row's price/volume data and the
# "fixed_ticker" (column I) value were mistakenly pulled from other tickers'
# rows. Re-point every data row (2-24) back to the BZ (KANZHUN LTD ADR) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$data = @(
    @(2, 32.89522147852465, 38.93419647216797, 44.14833216591025, 32.65955257200576, 398688012),
    @(3, 32.89522147852465, 38.93419647216797, 44.14833216591025, 32.65955257200576, 398688012),
    @(4, 32.89522147852465, 38.93419647216797, 44.14833216591025, 32.65955257200576, 398688012),
    @(5, 32.89522147852465, 38.93419647216797, 44.14833216591025, 32.65955257200576, 398688012),
    @(6, 32.89522147852465, 38.93419647216797, 44.14833216591025, 32.65955257200576, 398688012),
    @(7, 32.89522147852465, 38.93419647216797, 44.14833216591025, 32.65955257200576, 398688012),
    @(8, 38.78690271060533, 33.82807159423828, 38.78690271060533, 27.39631824885721, 398688012),
    @(9, 35.43846439285922, 34.46633529663086, 38.58069527717858, 31.49103733607324, 398688012),
    @(10, 34.59398789760181, 29.74317169189453, 36.82300774954026, 24.54867183302684, 398688012),
    @(11, 26.46346973472976, 23.10521125793457, 28.96743358345567, 17.5768491557172, 398688012),
    @(12, 25.29495176097673, 23.0855712890625, 28.46663972226955, 22.14290131379661, 398688012),
    @(13, 16.48688997071209, 10.73268032073975, 19.34435540745924, 9.569072788389915, 398688012),
    @(14, 20.62088396693322, 23.85149002075196, 25.78592465197201, 20.59142488701677, 398688012),
    @(15, 18.72572567859768, 18.15619659423828, 19.6487562197715, 15.72096891441984, 398688012),
    @(16, 15.347829593572, 18.3427677154541, 18.99870750127328, 14.04184033963433, 398688012),
    @(17, 14.91577271013049, 14.53281402587891, 16.60472160763389, 14.17931254939126, 398688012),
    @(18, 15.90809422893064, 13.81283283233643, 16.06697641223375, 12.48219312665409, 398688012),
    @(19, 17.58628788791799, 19.65176010131836, 20.14826740784151, 17.5366379148757, 398688012),
    @(20, 18.88713779251686, 13.48513793945312, 19.3141343655566, 13.18226818908657, 398688012),
    @(21, 17.49691775968724, 14.44836235046387, 19.36378433956678, 14.00647021772064, 398688012),
    @(22, 13.56457861032944, 14.30933952331543, 15.43144607918914, 12.76023735463156, 398688012),
    @(23, 18.77095399077979, 15.19312286376953, 19.14035514723584, 13.18723301036554, 398688012),
    @(24, 17.59621690479154, 18.82755470275879, 20.08371937935827, 16.93089708958082, 398688012)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value = $row[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value = $row[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value = $row[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value = $row[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "BZ"      # I: fixed_ticker
}
